$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 46 (this shifts rows 46:125 down to 47:126,
# copying formatting from the row above - matches the target diff where
# every row from the old 46 downward moves down by one position).
$ws.Rows("46:46").Insert()

# Populate the newly inserted row 46 with the new weekly price-report entry.
$ws.Range("A46").Value = 4
$ws.Range("B46").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C46").Value = "Los Lagos"
$ws.Range("D46").Value = 45002
$ws.Range("E46").Value = 10
$ws.Range("F46").Value = "Fruta"
$ws.Range("G46").Value = 100104
$ws.Range("H46").Value = "Frutos de pepita"
$ws.Range("I46").Value = 100104003
$ws.Range("J46").Value = "Membrillo"
$ws.Range("K46").Value = "Champion"
$ws.Range("L46").Value = "Primera"
$ws.Range("M46").Value = 300
$ws.Range("N46").Value = 17000
$ws.Range("O46").Value = 18000
$ws.Range("P46").Value = 17500
$ws.Range("Q46").Value = "$/caja 18 kilos empedrada"
$ws.Range("R46").Value = "Región de O'Higgins"
$ws.Range("S46").Value = 972
$ws.Range("T46").Value = 18
